# Updates the cryptocurrency price/volume table on the active worksheet
# to reflect the latest scrape, matching the commit:
# "Updated symbol list on Tue Jan 10 12:51:22 UTC 2023 with GitHub Actions"
#
# Only the Price (column D) and Volume(1h) (column E) cells that actually
# changed between scrapes are touched; everything else (Coin, Link, Data,
# Hora, and rows whose values are unavailable "--") is left untouched.
#
# Price/percentage figures are stored as plain text in this sheet (so that
# values like "--" / "--%" can appear for unlisted coins), so each target
# cell is forced to Text format before the write - otherwise Excel would
# auto-convert "274.82" to a number or "-1.49%" to a percentage value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "274.82";     E = "-1.49%" }
    @{ Row = 3;  D = "26.66";      E = "-2.61%" }
    @{ Row = 4;                    E = "0.94%"  }
    @{ Row = 5;  D = "0.06303";    E = "-0.53%" }
    @{ Row = 6;  D = "6.925";      E = "-0.12%" }
    @{ Row = 7;  D = "1.298";      E = "35.12%" }
    @{ Row = 8;  D = "0.8685";     E = "-1.34%" }
    @{ Row = 9;  D = "0.1531";     E = "3.48%"  }
    @{ Row = 10; D = "0.05016";    E = "-1.87%" }
    @{ Row = 11; D = "0.07424";    E = "1.99%"  }
    @{ Row = 12; D = "0.02900";    E = "-8.11%" }
    @{ Row = 13; D = "0.09063";    E = "-0.02%" }
    @{ Row = 14; D = "0.001568";   E = "-0.04%" }
    @{ Row = 15; D = "0.0006338";  E = "0.73%"  }
    @{ Row = 16; D = "0.005984";   E = "-1.27%" }
    @{ Row = 17;                   E = "0.05%"  }
    @{ Row = 18; D = "3.309";      E = "-2.18%" }
    @{ Row = 19; D = "2.284";      E = "-0.46%" }
    @{ Row = 20;                   E = "0.82%"  }
    @{ Row = 21; D = "0.1319";     E = "-1.64%" }
    @{ Row = 22; D = "3.910";      E = "1.60%"  }
    @{ Row = 23; D = "0.04380";    E = "1.69%"  }
    @{ Row = 24;                   E = "-0.66%" }
    @{ Row = 26; D = "0.0001202";  E = "0.12%"  }
    @{ Row = 27; D = "0.0001618";  E = "-4.37%" }
    @{ Row = 40; D = "0.04064";    E = "-0.32%" }
    @{ Row = 41; D = "0.007033";   E = "4.77%"  }
    @{ Row = 42; D = "0.1168";     E = "0.57%"  }
    @{ Row = 43; D = "0.002093";   E = "-3.37%" }
    @{ Row = 44;                   E = "-10.59%" }
    @{ Row = 45; D = "0.00005224"; E = "-0.26%" }
    @{ Row = 47; D = "1.486";      E = "-37.47%" }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $cell = $ws.Range("D$($u.Row)")
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $cell = $ws.Range("E$($u.Row)")
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
    }
}
